# Auto-generated edit script for cryptos.xlsx price/volume update
# Commit: Updated cryptos list on Sat Jul 27 05:33:28 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) text updates ---
# These look numeric, so force Text format first to preserve exact string
# formatting (leading/trailing zeros, "." thousands separators), then clear
# the format back to Normal so no stray style is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.981.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.258.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.825.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.999.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.250.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "382.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.515"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.836"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0285"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.45"
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume 1h %) updates ---
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E20").Value = "  +3.38%  "
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +7.26%  "
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("E36").Value = "  -6.89%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("E40").Value = "  +4.96%  "
$ws.Range("E41").Value = "  +7.12%  "
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("E51").Value = "  +3.43%  "

# --- Rows that were reordered (coin + link + price + volume moved together) ---
# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.38"
$ws.Range("D43").Style = "Normal"
# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.53"
$ws.Range("D44").Style = "Normal"
# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.646.46"
$ws.Range("D46").Style = "Normal"
# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "344.91"
$ws.Range("D47").Style = "Normal"
